$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds date-serial values. Every row from 2 to 283
# has the same value (46074) which must be bumped to 46075 (i.e. +1 day).
$startRow = 2
$endRow = 283

for ($row = $startRow; $row -le $endRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $cell.Value2 + 1
}
